# Israel Premier League workbook update
# The underlying source data had a set of rows whose entire record (all
# columns except the leading index column A) belongs to a different row
# than where it was previously placed. This script corrects that by
# rotating the B:AD column values among each affected group of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$groups = @(
    @(8,9),
    @(18,19),
    @(31,32),
    @(58,59),
    @(73,74),
    @(86,87),
    @(108,109,110),
    @(151,152),
    @(164,165),
    @(180,181),
    @(217,218),
    @(219,220,221),
    @(224,225,227),
    @(231,232),
    @(237,238),
    @(239,241)
)

foreach ($g in $groups) {
    $n = $g.Count
    $vals = @()
    foreach ($r in $g) {
        $rng = $ws.Range("B" + $r + ":AD" + $r)
        $vals += ,$rng.Value2
    }
    for ($i = 0; $i -lt $n; $i++) {
        $src = ($i + 1) % $n
        $destRow = $g[$i]
        $ws.Range("B" + $destRow + ":AD" + $destRow).Value2 = $vals[$src]
    }
}
